{"js": "// Apply CV wording/content tweaks by locating each original sentence/phrase\n// via Body.search() and replacing it in place (preserving the run's\n// formatting, since insertText(..., \"Replace\") replaces the matched range's\n// text without touching surrounding runs).\n\nconst replacements = [\n  {\n    find: \"A motivated undergraduate student looking to enrich the customer service with strong communication skills and a dedication to helping others.\",\n    replace: \"A motivated undergraduate student with strong communication skills and a dedication to helping others.\"\n  },\n  {\n    find: \"Demonstrated a passion for customer service, effectively communicating with students and parents to address weaker points in the student's education.\",\n    replace: \"Effectively communicated with students and parents to address weaker points in the student's education.\"\n  },\n  {\n    find: \"Developed public speaking skills when presenting findings to staff and actively participating in meetings.\",\n    replace: \"Applied and developed public speaking skills when presenting findings to staff and actively participating in meetings.\"\n  },\n  {\n    find: \"Participated in community service events and projects with other cadets, enhancing teamwork skills and the ability to positively interact with the community.\",\n    replace: \"Participated in community service events and projects with other cadets, enhancing teamwork skills and the ability to interact positively with the community.\"\n  },\n  {\n    find: \"in Mathematics, second year student.\",\n    replace: \"in Mathematics, second-year student: on track for First Class (97% average).\"\n  },\n  {\n    find: \"Developed organisational and teamwork skills collaborating on a group project, communicating effectively with other group members.\",\n    replace: \"Applied organisational and teamwork skills collaborating on a group project, communicating effectively with other group members.\"\n  },\n  {\n    find: \"A-Levels (AABCC)\",\n    replace: \"A-Levels (2 A's, 1 B, 2 C's)\"\n  },\n  {\n    find: \"Held multiple leadership roles, including Secretary of the Student Union managing communication between students and staff and organising the Union, as well as an Ambassador role to represent my college during open days and in the community.\",\n    replace: \"Held multiple leadership roles, including Secretary of the Student Union managing communication between students and staff and organising the Union, as well as an Ambassador role to represent my college during open days and within the community.\"\n  },\n  {\n    find: \"GCSEs (A88888887777)\",\n    replace: \"GCSEs (1 A, 7 8's, 4 7's)\"\n  },\n  {\n    find: \"Participated in foreign exchange programs, enhancing cultural awareness and the ability to adapt in new situations.\",\n    replace: \"Participated in foreign exchange programs, enhancing cultural awareness and the ability to adapt to new situations.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply CV wording/content tweaks by locating each original sentence/phrase\n# with Find.Execute and then overwriting just that matched Range's Text.\n# (Using Find.Execute's own Replace:=wdReplaceAll path would trigger Word's\n# AutoFormat \"smart quotes\" substitution on the straight apostrophes we are\n# inserting, so we set Range.Text directly on the found range instead - it\n# keeps the surrounding run formatting intact and leaves straight quotes\n# untouched.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-CvText($findText, $replaceText) {\n    $r = $d.Content\n    $find = $r.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n    $r.Text = $replaceText\n}\n\nReplace-CvText \"A motivated undergraduate student looking to enrich the customer service with strong communication skills and a dedication to helping others.\" \"A motivated undergraduate student with strong communication skills and a dedication to helping others.\"\n\nReplace-CvText \"Demonstrated a passion for customer service, effectively communicating with students and parents to address weaker points in the student's education.\" \"Effectively communicated with students and parents to address weaker points in the student's education.\"\n\nReplace-CvText \"Developed public speaking skills when presenting findings to staff and actively participating in meetings.\" \"Applied and developed public speaking skills when presenting findings to staff and actively participating in meetings.\"\n\nReplace-CvText \"Participated in community service events and projects with other cadets, enhancing teamwork skills and the ability to positively interact with the community.\" \"Participated in community service events and projects with other cadets, enhancing teamwork skills and the ability to interact positively with the community.\"\n\nReplace-CvText \"in Mathematics, second year student.\" \"in Mathematics, second-year student: on track for First Class (97% average).\"\n\nReplace-CvText \"Developed organisational and teamwork skills collaborating on a group project, communicating effectively with other group members.\" \"Applied organisational and teamwork skills collaborating on a group project, communicating effectively with other group members.\"\n\nReplace-CvText \"A-Levels (AABCC)\" \"A-Levels (2 A's, 1 B, 2 C's)\"\n\nReplace-CvText \"Held multiple leadership roles, including Secretary of the Student Union managing communication between students and staff and organising the Union, as well as an Ambassador role to represent my college during open days and in the community.\" \"Held multiple leadership roles, including Secretary of the Student Union managing communication between students and staff and organising the Union, as well as an Ambassador role to represent my college during open days and within the community.\"\n\nReplace-CvText \"GCSEs (A88888887777)\" \"GCSEs (1 A, 7 8's, 4 7's)\"\n\nReplace-CvText \"Participated in foreign exchange programs, enhancing cultural awareness and the ability to adapt in new situations.\" \"Participated in foreign exchange programs, enhancing cultural awareness and the ability to adapt to new situations.\"\n"}
